# Auto-generated edit script: update cryptos list price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '57.636.43'
$c.ClearFormats()
$ws.Range("E2").Value = '  -0.92%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.433.83'
$c.ClearFormats()
$ws.Range("E3").Value = '  -1.53%  '

$ws.Range("E4").Value = '  -0.03%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '505.88'
$c.ClearFormats()
$ws.Range("E5").Value = '  -2.79%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '128.89'
$c.ClearFormats()
$ws.Range("E6").Value = '  -2.80%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -1.34%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.448.36'
$c.ClearFormats()
$ws.Range("E9").Value = '  -1.00%  '

$ws.Range("E10").Value = '  -0.21%  '

$ws.Range("E11").Value = '  -3.96%  '

$ws.Range("E12").Value = '  -3.36%  '

$ws.Range("E13").Value = '  -3.33%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.868.33'
$c.ClearFormats()
$ws.Range("E14").Value = '  -1.40%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '57.571.54'
$c.ClearFormats()
$ws.Range("E15").Value = '  -0.92%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '21.76'
$c.ClearFormats()
$ws.Range("E16").Value = '  -1.57%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.0000132'
$c.ClearFormats()
$ws.Range("E17").Value = '  -2.87%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.444.54'
$c.ClearFormats()
$ws.Range("E18").Value = '  -1.26%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.43'
$c.ClearFormats()
$ws.Range("E19").Value = '  -3.86%  '

$ws.Range("E20").Value = '  -1.62%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '314.69'
$c.ClearFormats()
$ws.Range("E21").Value = '  -1.59%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.23%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.67'
$c.ClearFormats()
$ws.Range("E23").Value = '  -1.34%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '63.47'
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.58%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.408'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.22%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  -0.98%  '

$ws.Range("E28").Value = '  -2.43%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '169.91'
$c.ClearFormats()
$ws.Range("E29").Value = '  +2.81%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0₃0722'
$c.ClearFormats()
$ws.Range("E30").Value = '  -3.86%  '

$ws.Range("E31").Value = '  -3.00%  '

$ws.Range("E32").Value = '  -3.02%  '

$ws.Range("E33").Value = '  +0.19%  '

$ws.Range("E34").Value = '  -0.05%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.ClearFormats()
$ws.Range("E35").Value = '  -0.10%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '17.70'
$c.ClearFormats()
$ws.Range("E36").Value = '  -2.45%  '

$ws.Range("E37").Value = '  -5.31%  '

$ws.Range("E38").Value = '  -2.10%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '36.31'
$c.ClearFormats()
$ws.Range("E39").Value = '  -0.55%  '

$ws.Range("E40").Value = '  -2.64%  '

$ws.Range("E41").Value = '  -4.69%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '270.54'
$c.ClearFormats()
$ws.Range("E42").Value = '  -2.23%  '

$ws.Range("E43").Value = '  -2.88%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '4.88'
$c.ClearFormats()
$ws.Range("E44").Value = '  -3.33%  '

$ws.Range("E45").Value = '  -2.46%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0910'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.17%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '119.54'
$c.ClearFormats()
$ws.Range("E47").Value = '  -5.30%  '

$ws.Range("E48").Value = '  -1.63%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '17.15'
$c.ClearFormats()
$ws.Range("E49").Value = '  -3.87%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0209'
$c.ClearFormats()
$ws.Range("E50").Value = '  -2.27%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '16.59'
$c.ClearFormats()
$ws.Range("E51").Value = '  -3.27%  '
